$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 17 by copying row 16's formatting/types (numeric cells land
# with no explicit style, matching the rest of the data rows), then
# overwrite with the new values.
$ws.Range("A16:T16").Copy()
$ws.Range("A17:T17").PasteSpecial()

# Date/"week" text columns: Excel would normally auto-convert a
# "2024-01-04"-shaped string to a date serial (and "00" to the number 0)
# if assigned via .Value, which would also pull in a new number-format
# style. Copy the literal text from existing cells instead so the new
# row stays plain inline/shared text exactly like its neighbours.
$ws.Cells.Item(15, 1).Copy()
$ws.Cells.Item(17, 1).PasteSpecial()

$ws.Cells.Item(15, 4).Copy()
$ws.Cells.Item(17, 4).PasteSpecial()

$ws.Cells.Item(17, 2).Value = "18:44:22"
$ws.Cells.Item(17, 3).Value = "Thursday"

$ws.Cells.Item(17, 5).Value = 140502
$ws.Cells.Item(17, 6).Value = 142867
$ws.Cells.Item(17, 7).Value = 172271
$ws.Cells.Item(17, 8).Value = 147190
$ws.Cells.Item(17, 9).Value = -1
$ws.Cells.Item(17, 10).Value = 117911
$ws.Cells.Item(17, 11).Value = 224327
$ws.Cells.Item(17, 12).Value = 248466
$ws.Cells.Item(17, 13).Value = 184509
$ws.Cells.Item(17, 14).Value = 110002
$ws.Cells.Item(17, 15).Value = 40354
$ws.Cells.Item(17, 16).Value = 30804
$ws.Cells.Item(17, 17).Value = 72344
$ws.Cells.Item(17, 18).Value = -1
$ws.Cells.Item(17, 19).Value = 41606
$ws.Cells.Item(17, 20).Value = -1
